$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 28: "Robot instructions" (33.01)
$ws.Range("C28").Value = 'Recursive approach, base case ""'
$ws.Range("D28").Value = (Get-Date -Year 2025 -Month 7 -Day 31 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E28").Value = "10 minutes"
$ws.Range("F28").Value = "10 minutes"
$ws.Range("G28").Value = "Fine"
$ws.Range("H28").Value = "No"
$ws.Range("I28").Value = "That problem was breaking down into a subset of itself"
$ws.Range("J28").Value = "No"
$ws.Range("K28").Value = "No"
$ws.Range("L28").Value = "Good practice"
$ws.Range("M28").Value = "No"
$ws.Range("N28").Value = 4
$ws.Range("O28").Value = 4
$ws.Range("P28").Value = 4
$ws.Range("Q28").Value = 4

# Row 29: "Nested array sum" (33.02)
$ws.Range("C29").Value = "Used a helper function to keep track of state"
$ws.Range("D29").Value = (Get-Date -Year 2025 -Month 8 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E29").Value = "15 minutes"
$ws.Range("F29").Value = "15 minutes"
$ws.Range("G29").Value = "Fine"
$ws.Range("H29").Value = "Yes"
$ws.Range("I29").Value = "Eager/lazy checking"
$ws.Range("J29").Value = "My recursion was not very elegant"
$ws.Range("K29").Value = "No"
$ws.Range("L29").Value = "Good review"
$ws.Range("M29").Value = "Lazy > Eager for cases where there is more than one recursive call per element"
$ws.Range("N29").Value = 3
$ws.Range("O29").Value = 3
$ws.Range("P29").Value = 3
$ws.Range("Q29").Value = 3

# Row heights auto-adjusted by Excel after wrapped text entry
$ws.Rows.Item(28).RowHeight = 40
$ws.Rows.Item(29).RowHeight = 60

# Scroll / selection changes to match final view state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("B28").Select()
